# Refresh the cryptocurrency price/volume snapshot in the worksheet.
# (mirrors the automated "Updated cryptos list ... with GitHub Actions" commit)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '28.517.02'
$ws.Range('E2').Value = '  +0.86%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.873.22'
$ws.Range('E3').Value = '  +0.15%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.9994'
$ws.Range('E4').Value = '  -0.21%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '331.00'
$ws.Range('E5').Value = '  -1.69%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.9988'
$ws.Range('E6').Value = '  -0.22%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4589'
$ws.Range('E7').Value = '  -2.25%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.4019'
$ws.Range('E8').Value = '  +2.54%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '47.57'
$ws.Range('E9').Value = '  +0.49%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.07862'
$ws.Range('E10').Value = '  -1.52%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.9880'
$ws.Range('E11').Value = '  -1.63%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '21.43'
$ws.Range('E12').Value = '  -1.13%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '1.859.69'
$ws.Range('E13').Value = '  -1.31%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '5.841'
$ws.Range('E14').Value = '  -2.34%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '7.004'
$ws.Range('E15').Value = '  -3.62%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '1.000'
$ws.Range('E16').Value = '  -0.24%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '88.46'
$ws.Range('E17').Value = '  -2.99%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.06539'
$ws.Range('E18').Value = '  -0.69%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.00001018'
$ws.Range('E19').Value = '  -2.26%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '17.20'
$ws.Range('E20').Value = '  -2.65%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '0.9979'
$ws.Range('E21').Value = '  -0.25%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '28.489.93'
$ws.Range('E22').Value = '  +0.74%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '5.334'
$ws.Range('E23').Value = '  -1.97%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '10.86'
$ws.Range('E24').Value = '  -1.55%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.249'
$ws.Range('E25').Value = '  -2.08%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.083.91'
$ws.Range('E26').Value = '  -0.83%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '157.48'
$ws.Range('E27').Value = '  -0.96%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '19.30'
$ws.Range('E28').Value = '  -2.71%  '
$ws.Range('E29').Value = '  -3.47%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '5.313'
$ws.Range('E30').Value = '  -3.40%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '117.31'
$ws.Range('E31').Value = '  -2.07%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.9559'
$ws.Range('E32').Value = '  -2.10%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.09323'
$ws.Range('E33').Value = '  -1.60%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '3.598'
$ws.Range('E34').Value = '  +0.53%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.396'
$ws.Range('E35').Value = '  +1.91%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '5.228'
$ws.Range('E36').Value = '  -2.29%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.06015'
$ws.Range('E37').Value = '  -1.21%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.02205'
$ws.Range('E38').Value = '  -2.49%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '8.290'
$ws.Range('E39').Value = '  -1.59%  '
$ws.Range('E40').Value = '  -0.66%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.9978'
$ws.Range('E41').Value = '  -0.17%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.5766'
$ws.Range('E42').Value = '  -3.09%  '
$ws.Range('B43').Value = 'Algorand'
$ws.Range('C43').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.1810'
$ws.Range('E43').Value = '  -3.68%  '
$ws.Range('B44').Value = 'Aptos'
$ws.Range('C44').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '10.05'
$ws.Range('E44').Value = '  -3.10%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '1.254'
$ws.Range('E45').Value = '  -3.85%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '2.301'
$ws.Range('E46').Value = '  +16.32%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.5435'
$ws.Range('E47').Value = '  -3.08%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '11.83'
$ws.Range('E48').Value = '  -2.48%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.07190'
$ws.Range('E49').Value = '  +4.33%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.879'
$ws.Range('E50').Value = '  -4.41%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '109.94'
$ws.Range('E51').Value = '  -0.67%  '
